$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: new "saldo" (balance) log entry ---------------------------------
$ws.Range("A8").NumberFormat = "mm-dd-yy"
$ws.Range("A8").Value = 41761
$ws.Range("B8").Value = "Adicionada tabela com saldo"
$ws.Range("C8").Value = "não"

# Leave the selection where the user ended up after adding the new row --
$ws.Range("C12").Select() | Out-Null
